$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set column C width (stored OOXML width of 11)
$ws.Columns.Item(3).ColumnWidth = 10.1

# Header
$ws.Range("C1").Value = "Retenção IR"

# Values for C2:C17
$values = @("-", "-", "-", "-", "-", "-", "-", "-", "X", "-", "-", "-", "X", "-", "-", "X")

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 3)
    $cell.Value = $values[$i]
    $cell.HorizontalAlignment = -4108  # xlCenter
}

# Row 18: empty numeric-styled cells in B18 and C18
$ws.Cells.Item(18, 2).HorizontalAlignment = -4108
$ws.Cells.Item(18, 3).HorizontalAlignment = -4108
